$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
function tryprop($obj, $name) {
  try {
    $v = $obj.$name
    Write-Host "$name =" $v
  } catch {
    Write-Host "$name ERR:" $_
  }
}
tryprop $d "Name"
try {
  $d.Name = "Office Theme"
  Write-Host "after set Name:" $d.Name
} catch { Write-Host "set ERR:" $_ }
